$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1=14 and Q1=15, matching the bold/bordered
#     header style already used by O1 (and the rest of row 1). ---
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("P1").Value = 14

$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: swap I<->K values and M<->O values, then append
#     new columns P and Q (both = 2). ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2   # P (new) = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q (new) = 2
}
